$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 1.02
$ws.Cells.Item(2, 3).Value = 1.060212565978399
$ws.Cells.Item(2, 4).Value = 1.059029976307234
$ws.Cells.Item(2, 5).Value = 1.065305517444206
$ws.Cells.Item(2, 6).Value = 1.074369071146585
$ws.Cells.Item(2, 9).Value = 1.048266099844685
$ws.Cells.Item(2, 10).Value = 1.065194583162613
$ws.Cells.Item(2, 11).Value = 1.061760593960952
$ws.Cells.Item(2, 12).Value = 1.06801911496176
$ws.Cells.Item(2, 13).Value = 1.077058466466712
$ws.Cells.Item(3, 2).Value = 1.02
$ws.Cells.Item(3, 3).Value = 1.061496188969643
$ws.Cells.Item(3, 4).Value = 1.060029360536541
$ws.Cells.Item(3, 5).Value = 1.06648871049175
$ws.Cells.Item(3, 6).Value = 1.075724833772661
$ws.Cells.Item(3, 9).Value = 1.048655696046671
$ws.Cells.Item(3, 10).Value = 1.06613017084674
$ws.Cells.Item(3, 11).Value = 1.062573658481804
$ws.Cells.Item(3, 12).Value = 1.069016761748611
$ws.Cells.Item(3, 13).Value = 1.078230018876174
$ws.Cells.Item(4, 2).Value = 1.02
$ws.Cells.Item(4, 3).Value = 1.062326228587303
$ws.Cells.Item(4, 4).Value = 1.060675445267672
$ws.Cells.Item(4, 5).Value = 1.067254098604845
$ws.Cells.Item(4, 6).Value = 1.076602188036038
$ws.Cells.Item(4, 9).Value = 1.04890627515856
$ws.Cells.Item(4, 10).Value = 1.066734492371148
$ws.Cells.Item(4, 11).Value = 1.063098569188278
$ws.Cells.Item(4, 12).Value = 1.069661512171498
$ws.Cells.Item(4, 13).Value = 1.078987631889341
$ws.Cells.Item(5, 2).Value = 1.02
$ws.Cells.Item(5, 3).Value = 1.062675048634613
$ws.Cells.Item(5, 4).Value = 1.060946921260195
$ws.Cells.Item(5, 5).Value = 1.067575818131304
$ws.Cells.Item(5, 6).Value = 1.076971050995081
$ws.Cells.Item(5, 9).Value = 1.049011256702167
$ws.Cells.Item(5, 10).Value = 1.066988295935872
$ws.Cells.Item(5, 11).Value = 1.063318956924477
$ws.Cells.Item(5, 12).Value = 1.069932377336074
$ws.Cells.Item(5, 13).Value = 1.079306024841855
$ws.Cells.Item(6, 2).Value = 1.02
$ws.Cells.Item(6, 3).Value = 1.062733609604775
$ws.Cells.Item(6, 4).Value = 1.060992495229391
$ws.Cells.Item(6, 5).Value = 1.067629833411249
$ws.Cells.Item(6, 6).Value = 1.07703298615326
$ws.Cells.Item(6, 9).Value = 1.049028862369759
$ws.Cells.Item(6, 10).Value = 1.067030895883026
$ws.Cells.Item(6, 11).Value = 1.063355944347983
$ws.Cells.Item(6, 12).Value = 1.069977845831115
$ws.Cells.Item(6, 13).Value = 1.079359478184987
$ws.Cells.Item(7, 2).Value = 1.02
$ws.Cells.Item(7, 3).Value = 1.062330890039279
$ws.Cells.Item(7, 4).Value = 1.060679073283403
$ws.Cells.Item(7, 5).Value = 1.067258397629084
$ws.Cells.Item(7, 6).Value = 1.076607116706149
$ws.Cells.Item(7, 9).Value = 1.04890767934799
$ws.Cells.Item(7, 10).Value = 1.066737884699155
$ws.Cells.Item(7, 11).Value = 1.06310151513666
$ws.Cells.Item(7, 12).Value = 1.069665132220823
$ws.Cells.Item(7, 13).Value = 1.078991886691067
$ws.Cells.Item(8, 2).Value = 1.02
$ws.Cells.Item(8, 3).Value = 1.060646488206925
$ws.Cells.Item(8, 4).Value = 1.059367844774465
$ws.Cells.Item(8, 5).Value = 1.06570542936734
$ws.Cells.Item(8, 6).Value = 1.074827241541855
$ws.Cells.Item(8, 9).Value = 1.048398080105304
$ws.Cells.Item(8, 10).Value = 1.065510991421919
$ws.Cells.Item(8, 11).Value = 1.062035621669735
$ws.Cells.Item(8, 12).Value = 1.068356439991606
$ws.Cells.Item(8, 13).Value = 1.077454494788274
$ws.Cells.Item(9, 2).Value = 1.02
$ws.Cells.Item(9, 3).Value = 1.05767400184968
$ws.Cells.Item(9, 4).Value = 1.0570527451848
$ws.Cells.Item(9, 5).Value = 1.062967136605161
$ws.Cells.Item(9, 6).Value = 1.071691399359272
$ws.Cells.Item(9, 9).Value = 1.047488450379145
$ws.Cells.Item(9, 10).Value = 1.063340789633619
$ws.Cells.Item(9, 11).Value = 1.060148141737401
$ws.Cells.Item(9, 12).Value = 1.066044182788121
$ws.Cells.Item(9, 13).Value = 1.074741780369689
$ws.Cells.Item(10, 2).Value = 1.02
$ws.Cells.Item(10, 3).Value = 1.055689217650335
$ws.Cells.Item(10, 4).Value = 1.055506177875342
$ws.Cells.Item(10, 5).Value = 1.061140275568743
$ws.Cells.Item(10, 6).Value = 1.069601011056387
$ws.Cells.Item(10, 9).Value = 1.04687413670083
$ws.Cells.Item(10, 10).Value = 1.061888303823558
$ws.Cells.Item(10, 11).Value = 1.05888350085767
$ws.Cells.Item(10, 12).Value = 1.064498391807703
$ws.Cells.Item(10, 13).Value = 1.072930701826593
$ws.Cells.Item(11, 2).Value = 1.02
$ws.Cells.Item(11, 3).Value = 1.054828997794873
$ws.Cells.Item(11, 4).Value = 1.054835720577301
$ws.Cells.Item(11, 5).Value = 1.060348875237239
$ws.Cells.Item(11, 6).Value = 1.068695850170939
$ws.Cells.Item(11, 9).Value = 1.046606245906017
$ws.Cells.Item(11, 10).Value = 1.061257983178738
$ws.Cells.Item(11, 11).Value = 1.058334373306913
$ws.Cells.Item(11, 12).Value = 1.06382799995354
$ws.Cells.Item(11, 13).Value = 1.072145831480218
$ws.Cells.Item(12, 2).Value = 1.02
$ws.Cells.Item(12, 3).Value = 1.054509350421025
$ws.Cells.Item(12, 4).Value = 1.054586562907232
$ws.Cells.Item(12, 5).Value = 1.060054857256316
$ws.Cells.Item(12, 6).Value = 1.068359628171484
$ws.Cells.Item(12, 9).Value = 1.046506454288681
$ws.Cells.Item(12, 10).Value = 1.061023643320153
$ws.Cells.Item(12, 11).Value = 1.058130170833649
$ws.Cells.Item(12, 12).Value = 1.063578825462841
$ws.Cells.Item(12, 13).Value = 1.071854193119241
$ws.Cells.Item(13, 2).Value = 1.02
$ws.Cells.Item(13, 3).Value = 1.054577921516276
$ws.Cells.Item(13, 4).Value = 1.054640013539114
$ws.Cells.Item(13, 5).Value = 1.060117927722948
$ws.Cells.Item(13, 6).Value = 1.068431749213377
$ws.Cells.Item(13, 9).Value = 1.046527872848665
$ws.Cells.Item(13, 10).Value = 1.061073919593189
$ws.Cells.Item(13, 11).Value = 1.058173983494061
$ws.Cells.Item(13, 12).Value = 1.063632281571343
$ws.Cells.Item(13, 13).Value = 1.071916755230631
$ws.Cells.Item(14, 2).Value = 1.02
$ws.Cells.Item(14, 3).Value = 1.054802578183818
$ws.Cells.Item(14, 4).Value = 1.054815127587263
$ws.Cells.Item(14, 5).Value = 1.060324572792232
$ws.Cells.Item(14, 6).Value = 1.068668058077305
$ws.Cells.Item(14, 9).Value = 1.046598002921999
$ws.Cells.Item(14, 10).Value = 1.061238616886746
$ws.Cells.Item(14, 11).Value = 1.058317498613476
$ws.Cells.Item(14, 12).Value = 1.063807406414777
$ws.Cells.Item(14, 13).Value = 1.072121726673243
$ws.Cells.Item(15, 2).Value = 1.02
$ws.Cells.Item(15, 3).Value = 1.054940980098828
$ws.Cells.Item(15, 4).Value = 1.054923005165434
$ws.Cells.Item(15, 5).Value = 1.060451886042766
$ws.Cells.Item(15, 6).Value = 1.068813655004492
$ws.Cells.Item(15, 9).Value = 1.046641174565337
$ws.Cells.Item(15, 10).Value = 1.061340064324349
$ws.Cells.Item(15, 11).Value = 1.058405892199251
$ws.Cells.Item(15, 12).Value = 1.063915285180547
$ws.Cells.Item(15, 13).Value = 1.072248002656693
$ws.Cells.Item(16, 2).Value = 1.02
$ws.Cells.Item(16, 3).Value = 1.055746290518433
$ws.Cells.Item(16, 4).Value = 1.05555065721801
$ws.Cells.Item(16, 5).Value = 1.06119279041911
$ws.Cells.Item(16, 6).Value = 1.069661083107726
$ws.Cells.Item(16, 9).Value = 1.046891875825855
$ws.Cells.Item(16, 10).Value = 1.061930106727078
$ws.Cells.Item(16, 11).Value = 1.058919912225458
$ws.Cells.Item(16, 12).Value = 1.064542861064379
$ws.Cells.Item(16, 13).Value = 1.072982776878889
$ws.Cells.Item(17, 2).Value = 1.02
$ws.Cells.Item(17, 3).Value = 1.056251224787447
$ws.Cells.Item(17, 4).Value = 1.055944155415449
$ws.Cells.Item(17, 5).Value = 1.061657442376085
$ws.Cells.Item(17, 6).Value = 1.070192647433422
$ws.Cells.Item(17, 9).Value = 1.047048627498005
$ws.Cells.Item(17, 10).Value = 1.062299852131059
$ws.Cells.Item(17, 11).Value = 1.059241932290571
$ws.Cells.Item(17, 12).Value = 1.064936238682391
$ws.Cells.Item(17, 13).Value = 1.073443502006855
$ws.Cells.Item(18, 2).Value = 1.02
$ws.Cells.Item(18, 3).Value = 1.056545667946985
$ws.Cells.Item(18, 4).Value = 1.056173600820469
$ws.Cells.Item(18, 5).Value = 1.061928431714367
$ws.Cells.Item(18, 6).Value = 1.070502699292044
$ws.Cells.Item(18, 9).Value = 1.047139875928965
$ws.Cells.Item(18, 10).Value = 1.062515384947565
$ws.Cells.Item(18, 11).Value = 1.059429613704124
$ws.Cells.Item(18, 12).Value = 1.065165587666518
$ws.Cells.Item(18, 13).Value = 1.073712171582964
$ws.Cells.Item(19, 2).Value = 1.02
$ws.Cells.Item(19, 3).Value = 1.056646052605707
$ws.Cells.Item(19, 4).Value = 1.05625182304709
$ws.Cells.Item(19, 5).Value = 1.062020826348938
$ws.Cells.Item(19, 6).Value = 1.070608418966553
$ws.Cells.Item(19, 9).Value = 1.047170958424933
$ws.Cells.Item(19, 10).Value = 1.062588853517664
$ws.Cells.Item(19, 11).Value = 1.059493583274223
$ws.Cells.Item(19, 12).Value = 1.065243772604864
$ws.Cells.Item(19, 13).Value = 1.07380377027462
$ws.Cells.Item(20, 2).Value = 1.02
$ws.Cells.Item(20, 3).Value = 1.056197058057252
$ws.Cells.Item(20, 4).Value = 1.055901944603072
$ws.Cells.Item(20, 5).Value = 1.061607593187575
$ws.Cells.Item(20, 6).Value = 1.070135615703383
$ws.Cells.Item(20, 9).Value = 1.047031828379148
$ws.Cells.Item(20, 10).Value = 1.06226019576928
$ws.Cells.Item(20, 11).Value = 1.05920739785286
$ws.Cells.Item(20, 12).Value = 1.064894043508918
$ws.Cells.Item(20, 13).Value = 1.073394077163326
$ws.Cells.Item(21, 2).Value = 1.02
$ws.Cells.Item(21, 3).Value = 1.054736425833377
$ws.Cells.Item(21, 4).Value = 1.054763564194969
$ws.Cells.Item(21, 5).Value = 1.060263722557756
$ws.Cells.Item(21, 6).Value = 1.068598471172181
$ws.Cells.Item(21, 9).Value = 1.046577359242218
$ws.Cells.Item(21, 10).Value = 1.061190123476783
$ws.Cells.Item(21, 11).Value = 1.058275243411534
$ws.Cells.Item(21, 12).Value = 1.063755840982589
$ws.Cells.Item(21, 13).Value = 1.072061370549332
$ws.Cells.Item(22, 2).Value = 1.02
$ws.Cells.Item(22, 3).Value = 1.05381735235664
$ws.Cells.Item(22, 4).Value = 1.054047123963297
$ws.Cells.Item(22, 5).Value = 1.059418448459406
$ws.Cells.Item(22, 6).Value = 1.067631976998144
$ws.Cells.Item(22, 9).Value = 1.046289966629256
$ws.Cells.Item(22, 10).Value = 1.060516106145081
$ws.Cells.Item(22, 11).Value = 1.057687817230077
$ws.Cells.Item(22, 12).Value = 1.063039274403962
$ws.Cells.Item(22, 13).Value = 1.071222849984827
$ws.Cells.Item(23, 2).Value = 1.02
$ws.Cells.Item(23, 3).Value = 1.054304639788189
$ws.Cells.Item(23, 4).Value = 1.054426989152226
$ws.Cells.Item(23, 5).Value = 1.059866576455679
$ws.Cells.Item(23, 6).Value = 1.068144337932157
$ws.Cells.Item(23, 9).Value = 1.046442475697203
$ws.Cells.Item(23, 10).Value = 1.060873532059886
$ws.Cells.Item(23, 11).Value = 1.057999351075967
$ws.Cells.Item(23, 12).Value = 1.063419229302021
$ws.Cells.Item(23, 13).Value = 1.071667423074169
$ws.Cells.Item(24, 2).Value = 1.02
$ws.Cells.Item(24, 3).Value = 1.056221533907146
$ws.Cells.Item(24, 4).Value = 1.055921018085319
$ws.Cells.Item(24, 5).Value = 1.061630117996019
$ws.Cells.Item(24, 6).Value = 1.070161385889108
$ws.Cells.Item(24, 9).Value = 1.047039419740694
$ws.Cells.Item(24, 10).Value = 1.062278115184366
$ws.Cells.Item(24, 11).Value = 1.059223002934276
$ws.Cells.Item(24, 12).Value = 1.064913110005174
$ws.Cells.Item(24, 13).Value = 1.07341641031907
$ws.Cells.Item(25, 2).Value = 1.02
$ws.Cells.Item(25, 3).Value = 1.058442997324306
$ws.Cells.Item(25, 4).Value = 1.057651804118408
$ws.Cells.Item(25, 5).Value = 1.063675275024911
$ws.Cells.Item(25, 6).Value = 1.072502045977657
$ws.Cells.Item(25, 9).Value = 1.047724998360826
$ws.Cells.Item(25, 10).Value = 1.063902831083201
$ws.Cells.Item(25, 11).Value = 1.060637206640577
$ws.Cells.Item(25, 12).Value = 1.066642702244385
$ws.Cells.Item(25, 13).Value = 1.075443529022865
